# Update the "Análise" worksheet with new sample data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Análise")

# Update header row
$ws.Range("C1").Value = "K"
$ws.Range("D1").Value = "Ti"

# Clear the now-unused columns E:I (headers + data)
$ws.Range("E1:I5").Clear()

# Update data rows with new values
$ws.Range("B2").Value = 25925.248
$ws.Range("C2").Value = 20124.517
$ws.Range("D2").Value = 6412.352

$ws.Range("B3").Value = 26410
$ws.Range("C3").Value = 20010
$ws.Range("D3").Value = 4570

$ws.Range("B4").Value = 28877.324
$ws.Range("C4").Value = 20885.962
$ws.Range("D4").Value = 5288.021

$ws.Range("B5").Value = 28565.036
$ws.Range("C5").Value = 19099.573
$ws.Range("D5").Value = 7410.386
